$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value2 = 128
$ws.Range("F5").Value2 = 1885
$ws.Range("F6").Value2 = 144
$ws.Range("F7").Value2 = 3967
$ws.Range("F8").Value2 = 520
$ws.Range("F9").Value2 = 1036
$ws.Range("F10").Value2 = 1298
$ws.Range("F11").Value2 = 649
$ws.Range("F14").Value2 = 2157
$ws.Range("F16").Value2 = 645655
$ws.Range("F17").Value2 = 1593
$ws.Range("F18").Value2 = 462
$ws.Range("F19").Value2 = 1406
$ws.Range("F20").Value2 = 662
$ws.Range("F22").Value2 = 1238
$ws.Range("F23").Value2 = 2147
$ws.Range("F24").Value2 = 1096
$ws.Range("F25").Value2 = 2654
$ws.Range("F26").Value2 = 1522
$ws.Range("F27").Value2 = 743
$ws.Range("F28").Value2 = 1494
$ws.Range("F29").Value2 = 515
$ws.Range("F30").Value2 = 1066
$ws.Range("F31").Value2 = 1068
$ws.Range("F32").Value2 = 72
$ws.Range("F33").Value2 = 1991
$ws.Range("F34").Value2 = 1310
$ws.Range("F35").Value2 = 1188
$ws.Range("F36").Value2 = 1809
$ws.Range("F37").Value2 = 1124
$ws.Range("F41").Value2 = 2529
$ws.Range("F45").Value2 = 868
$ws.Range("F46").Value2 = 134
$ws.Range("F49").Value2 = 19

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value2 = 464
$ws.Range("F10").Value2 = 92
$ws.Range("F11").Value2 = 144260
$ws.Range("F12").Value2 = 144260
$ws.Range("F18").Value2 = 223
$ws.Range("F19").Value2 = 326
$ws.Range("F23").Value2 = 108
$ws.Range("F24").Value2 = 75
$ws.Range("F25").Value2 = 95
$ws.Range("F26").Value2 = 88
$ws.Range("F27").Value2 = 511
$ws.Range("F28").Value2 = 88
$ws.Range("F32").Value2 = 306
$ws.Range("F39").Value2 = 9
$ws.Range("F40").Value2 = 181

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value2 = 3105
$ws.Range("F6").Value2 = 231
$ws.Range("F8").Value2 = 812
$ws.Range("F9").Value2 = 1132
$ws.Range("F10").Value2 = 623
$ws.Range("F11").Value2 = 1568
$ws.Range("F12").Value2 = 469
$ws.Range("F13").Value2 = 39
$ws.Range("F14").Value2 = 1782

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value2 = 812
$ws.Range("F3").Value2 = 623
$ws.Range("F5").Value2 = 1568
$ws.Range("F6").Value2 = 469
$ws.Range("F7").Value2 = 128
$ws.Range("F8").Value2 = 1782
$ws.Range("F9").Value2 = 3967
$ws.Range("F11").Value2 = 520
$ws.Range("F12").Value2 = 1298
$ws.Range("F13").Value2 = 649
$ws.Range("F15").Value2 = 2157
$ws.Range("F18").Value2 = 645662
$ws.Range("F19").Value2 = 464
$ws.Range("F20").Value2 = 92
$ws.Range("F21").Value2 = 1593
$ws.Range("F22").Value2 = 144260
$ws.Range("F23").Value2 = 462
$ws.Range("F24").Value2 = 1406
$ws.Range("F25").Value2 = 662
$ws.Range("F27").Value2 = 1238
$ws.Range("F28").Value2 = 2147
$ws.Range("F29").Value2 = 1096
$ws.Range("F30").Value2 = 2654
$ws.Range("F31").Value2 = 1522
$ws.Range("F33").Value2 = 1494
$ws.Range("F35").Value2 = 515
$ws.Range("F36").Value2 = 108
$ws.Range("F37").Value2 = 1066
$ws.Range("F38").Value2 = 1068
$ws.Range("F39").Value2 = 75
$ws.Range("F40").Value2 = 72
$ws.Range("F41").Value2 = 1991
$ws.Range("F42").Value2 = 1310
$ws.Range("F43").Value2 = 1188
$ws.Range("F44").Value2 = 1810
$ws.Range("F45").Value2 = 1124
$ws.Range("F46").Value2 = 306
$ws.Range("F47").Value2 = 306
$ws.Range("F48").Value2 = 2529
$ws.Range("F51").Value2 = 134
